$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper data: the two brand-new e2e entries that must be inserted (in order)
# before the existing trailing "d61b9a3c" row on every sheet, plus the values
# needed to refresh that now-shifted trailing row.
# ---------------------------------------------------------------------------

# === Sheet "Overview" (columns A..G) =======================================
$ws1 = $wb.Worksheets.Item("Overview")

# Insert two fresh rows above the last data row (row 5), copying its
# formatting (hyperlink-style column B, date-style column G) along the way.
$ws1.Rows("5:5").Copy()
$ws1.Rows("5:5").Insert()
$ws1.Rows("5:5").Copy()
$ws1.Rows("5:5").Insert()

# Row 5: 1eadbff6
$ws1.Range("A5").Value = "1eadbff6-f6a7-4f04-88ee-ac3d2e4c5549.md"
$ws1.Range("B5").Value = "e2e\1eadbff6-f6a7-4f04-88ee-ac3d2e4c5549.md"
$ws1.Range("C5").Value = ".md"
$ws1.Range("D5").Value = ""
$ws1.Range("E5").Value = "Ready for handoff"
$ws1.Range("F5").Value = "Ready for handoff"
$ws1.Range("G5").Value = "2016-08-16 12:40:50"

# Row 6: 38fdfe28
$ws1.Range("A6").Value = "38fdfe28-5d2b-49c5-9a5c-d94814134381.md"
$ws1.Range("B6").Value = "e2e\38fdfe28-5d2b-49c5-9a5c-d94814134381.md"
$ws1.Range("C6").Value = ".md"
$ws1.Range("D6").Value = ""
$ws1.Range("E6").Value = "Ready for handoff"
$ws1.Range("F6").Value = "Ready for handoff"
$ws1.Range("G6").Value = "2016-08-16 12:40:50"

# Row 7: d61b9a3c (was row 5) - values unchanged, row merely shifted down.
$ws1.Range("A7").Value = "d61b9a3c-8605-4203-b454-fbea71ac7e9b.md"
$ws1.Range("B7").Value = "e2e\d61b9a3c-8605-4203-b454-fbea71ac7e9b.md"
$ws1.Range("C7").Value = ".md"
$ws1.Range("D7").Value = ""
$ws1.Range("E7").Value = "Ready for handoff"
$ws1.Range("F7").Value = "Ready for handoff"
$ws1.Range("G7").Value = "2016-08-16 12:38:33"

# Grow the table / autofilter to match the new extent.
$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:G7"))

# Rebuild the hyperlinks on column B from scratch (wholesale delete, since
# removing a single hyperlink isn't supported) in the final row order.
$ws1.Range("B2").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c05ef8eca7de03671bce4e5f41f3a535436b9ef8/e2e/a93ce76c-ffd7-4ab8-abbc-06d092475df3.md", "", "", "e2e\a93ce76c-ffd7-4ab8-abbc-06d092475df3.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/825f8250230ece24c2b52f4cef33355ffa2b8f35/e2e/1b34b1d9-c08a-4ed6-b90b-4612acc7419f.md", "", "", "e2e\1b34b1d9-c08a-4ed6-b90b-4612acc7419f.md")
$ws1.Hyperlinks.Add($ws1.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/825f8250230ece24c2b52f4cef33355ffa2b8f35/e2e/1c6ff0da-558e-47d1-8404-6c35d564d400.md", "", "", "e2e\1c6ff0da-558e-47d1-8404-6c35d564d400.md")
$ws1.Hyperlinks.Add($ws1.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bbefea96f4aeb183aa9a79a3ba4f6fd0a24b4f36/e2e/1eadbff6-f6a7-4f04-88ee-ac3d2e4c5549.md", "", "", "e2e\1eadbff6-f6a7-4f04-88ee-ac3d2e4c5549.md")
$ws1.Hyperlinks.Add($ws1.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e6b1c5f7b2f0a6c57a6e4f4d1a1c9d3e5f7a9b1/e2e/38fdfe28-5d2b-49c5-9a5c-d94814134381.md", "", "", "e2e\38fdfe28-5d2b-49c5-9a5c-d94814134381.md")
$ws1.Hyperlinks.Add($ws1.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/883296f898e33019b21c2f506f6d3f8fb35cc1a5/e2e/d61b9a3c-8605-4203-b454-fbea71ac7e9b.md", "", "", "e2e\d61b9a3c-8605-4203-b454-fbea71ac7e9b.md")

# === Sheet "zh-cn" (columns A..P) ===========================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows("5:5").Copy()
$ws2.Rows("5:5").Insert()
$ws2.Rows("5:5").Copy()
$ws2.Rows("5:5").Insert()

# Row 5: 1eadbff6
$ws2.Range("A5").Value = "1eadbff6-f6a7-4f04-88ee-ac3d2e4c5549.md"
$ws2.Range("B5").Value = ".md"
$ws2.Range("C5").Value = "Ready for handoff"
$ws2.Range("D5").Value = "e2e"
$ws2.Range("E5").Value = "ht"
$ws2.Range("F5").Value = "False"
$ws2.Range("G5").Value = "1eadbff6-f6a7-4f04-88ee-ac3d2e4c5549.3f3a37f011b7b76dfa16a0b56d89e762efed1299.zh-cn.xlf"
$ws2.Range("H5").Value = "2016-08-16 12:40:45"
$ws2.Range("I5").Value = ""
$ws2.Range("J5").Value = ""
$ws2.Range("K5").Value = "0001-01-01 00:00:00"
$ws2.Range("L5").Value = ""
$ws2.Range("M5").Value = "True"
$ws2.Range("N5").Value = ""
$ws2.Range("O5").Value = "False"
$ws2.Range("P5").Value = ""

# Row 6: 38fdfe28
$ws2.Range("A6").Value = "38fdfe28-5d2b-49c5-9a5c-d94814134381.md"
$ws2.Range("B6").Value = ".md"
$ws2.Range("C6").Value = "Ready for handoff"
$ws2.Range("D6").Value = "e2e"
$ws2.Range("E6").Value = "ht"
$ws2.Range("F6").Value = "False"
$ws2.Range("G6").Value = "38fdfe28-5d2b-49c5-9a5c-d94814134381.94acb4bb0f912a43bc9e82f750e3724f799ec73a.zh-cn.xlf"
$ws2.Range("H6").Value = "2016-08-16 12:40:45"
$ws2.Range("I6").Value = ""
$ws2.Range("J6").Value = ""
$ws2.Range("K6").Value = "0001-01-01 00:00:00"
$ws2.Range("L6").Value = ""
$ws2.Range("M6").Value = "True"
$ws2.Range("N6").Value = ""
$ws2.Range("O6").Value = "False"
$ws2.Range("P6").Value = ""

# Row 7: d61b9a3c (was row 5) - values unchanged, row merely shifted down.
$ws2.Range("A7").Value = "d61b9a3c-8605-4203-b454-fbea71ac7e9b.md"
$ws2.Range("B7").Value = ".md"
$ws2.Range("C7").Value = "Ready for handoff"
$ws2.Range("D7").Value = "e2e"
$ws2.Range("E7").Value = "ht"
$ws2.Range("F7").Value = "False"
$ws2.Range("G7").Value = "d61b9a3c-8605-4203-b454-fbea71ac7e9b.1d24b9a3cc011368bacac029e6d8765221cd959c.zh-cn.xlf"
$ws2.Range("H7").Value = "2016-08-16 12:38:27"
$ws2.Range("I7").Value = ""
$ws2.Range("J7").Value = ""
$ws2.Range("K7").Value = "0001-01-01 00:00:00"
$ws2.Range("L7").Value = ""
$ws2.Range("M7").Value = "True"
$ws2.Range("N7").Value = ""
$ws2.Range("O7").Value = "False"
$ws2.Range("P7").Value = ""

$lo2 = $ws2.ListObjects.Item(1)
$lo2.Resize($ws2.Range("A1:P7"))

$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c05ef8eca7de03671bce4e5f41f3a535436b9ef8/e2e/a93ce76c-ffd7-4ab8-abbc-06d092475df3.md", "", "", "a93ce76c-ffd7-4ab8-abbc-06d092475df3.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/683e7b795bed267c909466e8a10699b424edfbd6/e2e/a93ce76c-ffd7-4ab8-abbc-06d092475df3.md", "", "", "a93ce76c-ffd7-4ab8-abbc-06d092475df3.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/825f8250230ece24c2b52f4cef33355ffa2b8f35/e2e/1b34b1d9-c08a-4ed6-b90b-4612acc7419f.md", "", "", "1b34b1d9-c08a-4ed6-b90b-4612acc7419f.md")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/825f8250230ece24c2b52f4cef33355ffa2b8f35/e2e/1c6ff0da-558e-47d1-8404-6c35d564d400.md", "", "", "1c6ff0da-558e-47d1-8404-6c35d564d400.md")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bbefea96f4aeb183aa9a79a3ba4f6fd0a24b4f36/e2e/1eadbff6-f6a7-4f04-88ee-ac3d2e4c5549.md", "", "", "1eadbff6-f6a7-4f04-88ee-ac3d2e4c5549.md")
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e6b1c5f7b2f0a6c57a6e4f4d1a1c9d3e5f7a9b1/e2e/38fdfe28-5d2b-49c5-9a5c-d94814134381.md", "", "", "38fdfe28-5d2b-49c5-9a5c-d94814134381.md")
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/883296f898e33019b21c2f506f6d3f8fb35cc1a5/e2e/d61b9a3c-8605-4203-b454-fbea71ac7e9b.md", "", "", "d61b9a3c-8605-4203-b454-fbea71ac7e9b.md")

# === Sheet "de-de" (columns A..P) ===========================================
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows("5:5").Copy()
$ws3.Rows("5:5").Insert()
$ws3.Rows("5:5").Copy()
$ws3.Rows("5:5").Insert()

# Row 5: 1eadbff6
$ws3.Range("A5").Value = "1eadbff6-f6a7-4f04-88ee-ac3d2e4c5549.md"
$ws3.Range("B5").Value = ".md"
$ws3.Range("C5").Value = "Ready for handoff"
$ws3.Range("D5").Value = "e2e"
$ws3.Range("E5").Value = "ht"
$ws3.Range("F5").Value = "False"
$ws3.Range("G5").Value = "1eadbff6-f6a7-4f04-88ee-ac3d2e4c5549.3f3a37f011b7b76dfa16a0b56d89e762efed1299.de-de.xlf"
$ws3.Range("H5").Value = "2016-08-16 12:40:50"
$ws3.Range("I5").Value = ""
$ws3.Range("J5").Value = ""
$ws3.Range("K5").Value = "0001-01-01 00:00:00"
$ws3.Range("L5").Value = ""
$ws3.Range("M5").Value = "True"
$ws3.Range("N5").Value = ""
$ws3.Range("O5").Value = "False"
$ws3.Range("P5").Value = ""

# Row 6: 38fdfe28
$ws3.Range("A6").Value = "38fdfe28-5d2b-49c5-9a5c-d94814134381.md"
$ws3.Range("B6").Value = ".md"
$ws3.Range("C6").Value = "Ready for handoff"
$ws3.Range("D6").Value = "e2e"
$ws3.Range("E6").Value = "ht"
$ws3.Range("F6").Value = "False"
$ws3.Range("G6").Value = "38fdfe28-5d2b-49c5-9a5c-d94814134381.94acb4bb0f912a43bc9e82f750e3724f799ec73a.de-de.xlf"
$ws3.Range("H6").Value = "2016-08-16 12:40:50"
$ws3.Range("I6").Value = ""
$ws3.Range("J6").Value = ""
$ws3.Range("K6").Value = "0001-01-01 00:00:00"
$ws3.Range("L6").Value = ""
$ws3.Range("M6").Value = "True"
$ws3.Range("N6").Value = ""
$ws3.Range("O6").Value = "False"
$ws3.Range("P6").Value = ""

# Row 7: d61b9a3c (was row 5) - values unchanged, row merely shifted down.
$ws3.Range("A7").Value = "d61b9a3c-8605-4203-b454-fbea71ac7e9b.md"
$ws3.Range("B7").Value = ".md"
$ws3.Range("C7").Value = "Ready for handoff"
$ws3.Range("D7").Value = "e2e"
$ws3.Range("E7").Value = "ht"
$ws3.Range("F7").Value = "False"
$ws3.Range("G7").Value = "d61b9a3c-8605-4203-b454-fbea71ac7e9b.1d24b9a3cc011368bacac029e6d8765221cd959c.de-de.xlf"
$ws3.Range("H7").Value = "2016-08-16 12:38:33"
$ws3.Range("I7").Value = ""
$ws3.Range("J7").Value = ""
$ws3.Range("K7").Value = "0001-01-01 00:00:00"
$ws3.Range("L7").Value = ""
$ws3.Range("M7").Value = "True"
$ws3.Range("N7").Value = ""
$ws3.Range("O7").Value = "False"
$ws3.Range("P7").Value = ""

$lo3 = $ws3.ListObjects.Item(1)
$lo3.Resize($ws3.Range("A1:P7"))

$ws3.Range("A2").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c05ef8eca7de03671bce4e5f41f3a535436b9ef8/e2e/a93ce76c-ffd7-4ab8-abbc-06d092475df3.md", "", "", "a93ce76c-ffd7-4ab8-abbc-06d092475df3.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4689d3a4655107688241d3506400a1cee0ad677f/e2e/a93ce76c-ffd7-4ab8-abbc-06d092475df3.md", "", "", "a93ce76c-ffd7-4ab8-abbc-06d092475df3.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/825f8250230ece24c2b52f4cef33355ffa2b8f35/e2e/1b34b1d9-c08a-4ed6-b90b-4612acc7419f.md", "", "", "1b34b1d9-c08a-4ed6-b90b-4612acc7419f.md")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/825f8250230ece24c2b52f4cef33355ffa2b8f35/e2e/1c6ff0da-558e-47d1-8404-6c35d564d400.md", "", "", "1c6ff0da-558e-47d1-8404-6c35d564d400.md")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bbefea96f4aeb183aa9a79a3ba4f6fd0a24b4f36/e2e/1eadbff6-f6a7-4f04-88ee-ac3d2e4c5549.md", "", "", "1eadbff6-f6a7-4f04-88ee-ac3d2e4c5549.md")
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e6b1c5f7b2f0a6c57a6e4f4d1a1c9d3e5f7a9b1/e2e/38fdfe28-5d2b-49c5-9a5c-d94814134381.md", "", "", "38fdfe28-5d2b-49c5-9a5c-d94814134381.md")
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/883296f898e33019b21c2f506f6d3f8fb35cc1a5/e2e/d61b9a3c-8605-4203-b454-fbea71ac7e9b.md", "", "", "d61b9a3c-8605-4203-b454-fbea71ac7e9b.md")
